$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")
$ws.Activate()

# New key/value row for the SQL Server JDBC driver
$ws.Range("A8").Value = "dbdriver"
$ws.Range("B8").Value = "com.microsoft.sqlserver.jdbc.SQLServerDriver"

# Widen column B to fit the new long driver class name, keep C:D as before
$ws.Columns.Item(2).ColumnWidth = 41.285714285714285

# Update the active selection to reflect where the user last clicked
$ws.Range("B3").Select()
